$d = $word.ActiveDocument

# Locate the unique anchor paragraph - the list item reading:
#   "Run PRC_CONFIGURATION manually with input parameter as 1(for TASC)."
# The new bullet must be inserted immediately after it.
$target = "Run PRC_CONFIGURATION manually with input parameter as 1(for TASC)."
$newBulletMarker = "EDU_USER_ROLE_UPDATE.sql"
$count = $d.Paragraphs.Count
$anchorIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text.TrimEnd()
    if ($t -eq $target) {
        $anchorIndex = $i
    }
    if ($t.StartsWith($newBulletMarker)) {
        # Already applied - nothing further to do.
        Write-Output "EDU_USER_ROLE_UPDATE.sql bullet already present; no changes made."
        return
    }
}

if ($anchorIndex -eq -1) {
    throw "Anchor paragraph not found"
}

$anchorParagraph = $d.Paragraphs.Item($anchorIndex)

# Insert a brand-new paragraph right after the anchor paragraph. The new
# paragraph automatically inherits the anchor's pPr (ListParagraph style,
# numId 4, sz/szCs 20) plus a matching empty run.
$anchorParagraph.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($anchorIndex + 1)
$newRange = $newPara.Range
$newRange.Text = "EDU_USER_ROLE_UPDATE.sql"

# Append the remainder of the sentence onto the same paragraph (just before
# its paragraph mark).
$newPara2 = $d.Paragraphs.Item($anchorIndex + 1)
$tailStart = $newPara2.Range.End - 1
$tailRange = $d.Range($tailStart, $tailStart)
$tailRange.InsertAfter(" – Edu user role update for edu admin.")

Write-Output "Inserted EDU_USER_ROLE_UPDATE.sql bullet after TASC PRC_CONFIGURATION line (index $anchorIndex)."
